# Apply updated F-column (想去人数 / want-to-go count) values across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1269
$ws.Range("F5").Value = 1425
$ws.Range("F7").Value = 41
$ws.Range("F8").Value = 680
$ws.Range("F9").Value = 157
$ws.Range("F10").Value = 172
$ws.Range("F12").Value = 1311
$ws.Range("F14").Value = 532
$ws.Range("F15").Value = 477
$ws.Range("F16").Value = 130
$ws.Range("F18").Value = 125
$ws.Range("F19").Value = 756
$ws.Range("F20").Value = 2564
$ws.Range("F24").Value = 299
$ws.Range("F25").Value = 182
$ws.Range("F26").Value = 12
$ws.Range("F27").Value = 122
$ws.Range("F28").Value = 571
$ws.Range("F30").Value = 48
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 245

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 218
$ws.Range("F9").Value = 287
$ws.Range("F13").Value = 513
$ws.Range("F14").Value = 78
$ws.Range("F18").Value = 29

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2327
$ws.Range("F6").Value = 937
$ws.Range("F9").Value = 1176
$ws.Range("F10").Value = 295
$ws.Range("F11").Value = 83

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2327
$ws.Range("F9").Value = 937
$ws.Range("F10").Value = 1176
$ws.Range("F11").Value = 295
$ws.Range("F12").Value = 83
$ws.Range("F13").Value = 1269
$ws.Range("F14").Value = 1425
$ws.Range("F16").Value = 41
$ws.Range("F17").Value = 680
$ws.Range("F18").Value = 157
$ws.Range("F20").Value = 172
$ws.Range("F23").Value = 532
$ws.Range("F24").Value = 477
$ws.Range("F25").Value = 130
$ws.Range("F26").Value = 125
$ws.Range("F27").Value = 756
$ws.Range("F28").Value = 2564
$ws.Range("F31").Value = 299
$ws.Range("F32").Value = 287
$ws.Range("F33").Value = 182
$ws.Range("F34").Value = 122
$ws.Range("F35").Value = 571
$ws.Range("F37").Value = 513
$ws.Range("F38").Value = 78
$ws.Range("F39").Value = 48
$ws.Range("F41").Value = 29
$ws.Range("F48").Value = 32
$ws.Range("F49").Value = 245
